# Daily data update: add 3 new days of data at the top of the "June" sheet
# and add running-total helper formulas (last-year comparison) in columns J:M.

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("June")

# --- Insert 3 new rows above row 2, pushing all existing data down ---
$ws.Range("A2:A4").EntireRow.Insert()

# Copy the date number-format from the (now shifted) row 5 onto the new rows
$ws.Range("A5").Copy()
$ws.Range("A2:A4").PasteSpecial(-4122)  # xlPasteFormats
$excel.CutCopyMode = $false

# --- Fill in the 3 new rows of daily data ---
$ws.Range("A2").Value = 45460
$ws.Range("B2").Value = 259
$ws.Range("C2").Value = 26
$ws.Range("D2").Value = 7
$ws.Range("E2").Value = 226

$ws.Range("A3").Value = 45459
$ws.Range("B3").Value = 260
$ws.Range("C3").Value = 27
$ws.Range("D3").Value = 3
$ws.Range("E3").Value = 230

$ws.Range("A4").Value = 45458
$ws.Range("B4").Value = 221
$ws.Range("C4").Value = 28
$ws.Range("D4").Value = 7
$ws.Range("E4").Value = 186

# --- Running-total / prior-year comparison formulas (now on rows 6-8) ---
$ws.Range("J6").Formula = "=SUM(B2,B3,B4)"
$ws.Range("K6:M6").Formula = "=SUM(C2,C3,C4)"

$ws.Range("J7").Value = 3079
$ws.Range("K7").Value = 336
$ws.Range("L7").Value = 66

$ws.Range("J8").Formula = "=SUM(J7,J6)"
$ws.Range("K8").Formula = "=SUM(K7,K6)"
$ws.Range("L8").Formula = "=SUM(L7,L6)"

# The old "last year" helper formula that used to live in J6 shifted down to
# J9 when the 3 rows were inserted above; it has been superseded by the new
# J6:M8 block above, so clear it out.
$ws.Range("J9").ClearContents()

# --- Restore selection as last saved ---
$ws.Range("J10").Select()
